$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 78 (pushes existing rows 78..94 down to 79..95)
$ws.Rows("78:78").Insert()

# Populate the newly inserted row 78 with the new data record
$ws.Cells.Item(78, 1).Value = 8
$ws.Cells.Item(78, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(78, 3).Value = "Coquimbo"
$ws.Cells.Item(78, 4).Value = 44505
$ws.Cells.Item(78, 5).Value = 4
$ws.Cells.Item(78, 6).Value = 100112001
$ws.Cells.Item(78, 7).Value = "Berenjena"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 600
$ws.Cells.Item(78, 11).Value = 8000
$ws.Cells.Item(78, 12).Value = 9000
$ws.Cells.Item(78, 13).Value = 8500
$ws.Cells.Item(78, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(78, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(78, 16).Value = 142
$ws.Cells.Item(78, 17).Value = 60
$ws.Cells.Item(78, 18).Value = "Hortaliza"
